$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert two new worksheets ("Dataset_1" and "Dataset_2") right after
#    "Model_data_1" and before "Model_data_2".
# ---------------------------------------------------------------------------
$modelData1 = $wb.Worksheets.Item("Model_data_1")
$ds1 = $wb.Worksheets.Add([System.Type]::Missing, $modelData1)
$ds1.Name = "Dataset_1"
$ds2 = $wb.Worksheets.Add([System.Type]::Missing, $ds1)
$ds2.Name = "Dataset_2"

# ---------------------------------------------------------------------------
# 2) Populate "Dataset_1".
#    The order in which brand-new text values are first written matters: it
#    controls the order new entries land in the shared-strings table, which
#    needs to match the target file, so the two sheets' writes are
#    interleaved in a specific sequence below.
# ---------------------------------------------------------------------------
$ds1.Range("C2").Value = "Number of observations"
$ds1.Range("D2").Value = "Description step"
$ds1.Range("E2").Value = "Change"
$ds1.Range("D3").Value = "HRS dataset 2016"
$ds1.Range("D4").Value = "Filled in question about discrimination and measured health outcomes"
$ds1.Range("F2").Value = "Details"
$ds1.Range("D5").Value = "Remove unrealistic health outcome measurements"
$ds1.Range("F5").Value = "SYSBPM's <30, >500, wLBS of <30, >500, hIn of <40, >95, waist of <15, >95"

# "Dataset_2" row 5 description differs from "Dataset_1" -- introduce it now
# so it lands at the correct shared-string slot, ahead of the remaining
# strings common to both sheets ("Remove observations with missing values of
# moderators" / "Sex, wealth bin, race and age").
$ds2.Range("D5").Value = "Remove unrealistic health outcome measurements and not available differences"

$ds1.Range("D6").Value = "Remove observations with missing values of moderators"
$ds1.Range("F6").Value = "Sex, wealth bin, race and age"

# Header styling (bold font + thin bottom border) -- matches the style used
# for the other section headers in this workbook.
$hdr1 = $ds1.Range("C2:F2")
$hdr1.Font.Bold = $true
$hdr1.Borders.Item(9).LineStyle = 1
$hdr1.Borders.Item(9).Weight = 2

# Numeric data + formulas for "Dataset_1"
$ds1.Range("C3").Value = 20912
$ds1.Range("C4").Value = 5551
$ds1.Range("E4").Formula = "=C4-C3"
$ds1.Range("C5").Formula = "=C4+E5"
$ds1.Range("E5").Value = -39
$ds1.Range("C6").Formula = "=C5+E6"
$ds1.Range("E6").Value = -931

# Column widths (best-effort match of the bestFit widths from the source)
$ds1.Columns.Item(3).ColumnWidth = 22.0
$ds1.Columns.Item(4).ColumnWidth = 64.5
$ds1.Columns.Item(5).ColumnWidth = 6.67

$ds1.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$ds1.Range("F6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Populate "Dataset_2" (same layout, slightly different values/formulas).
# ---------------------------------------------------------------------------
$ds2.Range("C2").Value = "Number of observations"
$ds2.Range("D2").Value = "Description step"
$ds2.Range("E2").Value = "Change"
$ds2.Range("F2").Value = "Details"
$ds2.Range("D3").Value = "HRS dataset 2016"
$ds2.Range("D4").Value = "Filled in question about discrimination and measured health outcomes"
$ds2.Range("F5").Value = "SYSBPM's <30, >500, wLBS of <30, >500, hIn of <40, >95, waist of <15, >95"
$ds2.Range("D6").Value = "Remove observations with missing values of moderators"
$ds2.Range("F6").Value = "Sex, wealth bin, race and age"

$hdr2 = $ds2.Range("C2:F2")
$hdr2.Font.Bold = $true
$hdr2.Borders.Item(9).LineStyle = 1
$hdr2.Borders.Item(9).Weight = 2

$ds2.Range("C3").Value = 20912
$ds2.Range("C4").Value = 5551
$ds2.Range("E4").Formula = "=C4-C3"
$ds2.Range("C5").Formula = "=C4+E5"
$ds2.Range("E5").Value = -1354
$ds2.Range("C6").Value = 4103
$ds2.Range("E6").Formula = "=C6-C5"

$ds2.Columns.Item(3).ColumnWidth = 22.0
$ds2.Columns.Item(4).ColumnWidth = 64.5
$ds2.Columns.Item(5).ColumnWidth = 6.67

# ---------------------------------------------------------------------------
# 4) "Dataset_2" becomes the active sheet/tab (matches activeTab=3); this
#    also naturally clears any tabSelected flag previously on "Model_data_1".
# ---------------------------------------------------------------------------
$ds2.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$ds2.Range("H12").Select() | Out-Null

Write-Output "done"
